# Applies the motilal_portfolio_change_engine update:
# Inserts a new "Industry" column (C) between "Stock Name" (B) and "Mutual Fund"
# (which shifts from C to D, with all subsequent columns shifting right by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C ("Mutual Fund"), shifting
# Mutual Fund/Status/Jan_2026/Dec_2025/Oct_2025/MoM/QoQ each one column to the right.
$ws.Columns.Item(3).Insert()

# Set the new header
$ws.Range("C1").Value = "Industry"

# Populate the Industry values for each holding row (2-23)
$industries = @{
    2  = "Banks"
    3  = "Banks"
    4  = "Power"
    5  = "Telecom - Services"
    6  = "Pharmaceuticals & Biotechnology"
    7  = "Finance"
    8  = "Petroleum Products"
    9  = "Healthcare Services"
    10 = "Finance"
    11 = "Consumer Durables"
    12 = "Retailing"
    13 = "Minerals & Mining"
    14 = "Banks"
    15 = "Ferrous Metals"
    16 = "Finance"
    17 = "Electrical Equipment"
    18 = "Banks"
    19 = "Healthcare Services"
    20 = "IT - Software"
    21 = "Realty"
    22 = "Chemicals & Petrochemicals"
    23 = "Automobiles"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
